# Apply data updates to the ESIS.table sheet (added poison sumac to shrubs list,
# which shifted/recalculated several min/max percentile values throughout the table)
# and restore the sheet view/selection to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("ESIS.table")

# --- Cell value updates -----------------------------------------------------

$updates = @{
    "F17" = 5
    "N17" = 0.1
    "O17" = 0.2
    "Q17" = 5

    "N18" = 0.1
    "O18" = 0.2

    "F19" = 45
    "N19" = 0.1
    "O19" = 0.2
    "Q19" = 45

    "C20" = 2
    "D20" = 10
    "F20" = 45
    "N20" = 0.5
    "O20" = 1.5
    "Q20" = 45

    "C21" = 20
    "F21" = 40
    "N21" = 3
    "O21" = 10
    "Q21" = 40

    "F29" = 0.4
    "P29" = 0.5
    "Q29" = 0.5

    "P30" = 0
    "Q30" = 0

    "D38" = 0
    "O38" = 0
    "P38" = 0.5
    "Q38" = 5

    "D39" = 0
    "O39" = 0
    "P39" = 0.5
    "Q39" = 5

    "D40" = 0
    "E40" = 0.1
    "O40" = 0
    "P40" = 1.5

    "D41" = 0.5
    "O41" = 1
    "P41" = 1.5

    "D42" = 0.5
    "O42" = 1
    "P42" = 1.5

    "F52" = 0.2
    "O52" = 4

    "F53" = 0.2
    "O53" = 4
    "Q53" = 0.4

    "O54" = 4

    "O55" = 0.5
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Restore view / selection -----------------------------------------------

$ws.Activate()
$ws.Range("C3").Select()
